# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "23.58") are stored as text, matching the original inlineStr cells,
# then restore the style so cells do not carry a lingering format index.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '60.916.93'
$ws.Range('E2').Value = '  +2.83%  '
$ws.Range('D3').Value = '2.607.25'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('D5').Value = '571.02'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '143.21'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').Value = '2.633.39'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('D11').Value = '0.106'
$ws.Range('E11').Value = '  +2.79%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  +6.87%  '
$ws.Range('D14').Value = '3.072.75'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').Value = '60.909.19'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = '23.58'
$ws.Range('E16').Value = '  +5.06%  '
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '2.621.90'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('D19').Value = '11.24'
$ws.Range('E19').Value = '  +9.55%  '
$ws.Range('D20').Value = '4.65'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').Value = '349.15'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '7.09'
$ws.Range('E22').Value = '  +12.76%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '0.522'
$ws.Range('E24').Value = '  +12.87%  '
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +5.87%  '
$ws.Range('D29').Value = '0.0₃0798'
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').Value = '1.84'
$ws.Range('E30').Value = '  +8.77%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  +4.64%  '
$ws.Range('D33').Value = '159.95'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').Value = '19.52'
$ws.Range('E34').Value = '  +2.65%  '
$ws.Range('E35').Value = '  +5.73%  '
$ws.Range('D36').Value = '0.968'
$ws.Range('E36').Value = '  +10.26%  '
$ws.Range('E37').Value = '  +3.78%  '
$ws.Range('E38').Value = '  +7.18%  '
$ws.Range('D39').Value = '37.83'
$ws.Range('D40').Value = '0.854'
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('E41').Value = '  +3.73%  '
$ws.Range('D42').Value = '298.78'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').Value = '140.36'
$ws.Range('E43').Value = '  +7.00%  '
$ws.Range('D44').Value = '0.0989'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '0.607'
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').Value = '0.0549'
$ws.Range('E47').Value = '  +2.13%  '
$ws.Range('D48').Value = '0.0242'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '19.55'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '4.94'
$ws.Range('E50').Value = '  +9.19%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '10.71'
$ws.Range('E51').Value = '  +0.51%  '

# Restore normal style on column D (clears the temporary text format marker)
$ws.Range("D2:D51").Style = "Normal"
